$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = 815
